# docs/课程安排.xlsx — feat(widget): add common properties
#
# 1. Sheet view: zoom to 90%, reset the frozen-pane scroll anchor back to
#    A2 (top of the data, just below the frozen header row), and move the
#    active selection in the lower (frozen-below) pane to E15.
# 2. Column layout: column D-J used to be one uniform 20.1796875-wide
#    block (D:I) plus a slightly wider J. Re-slice that block into narrow
#    "index" columns (D, F, G, I) flanking two wide text columns (E, H),
#    and widen J a bit.

$wb  = $excel.ActiveWorkbook
$ws  = $wb.ActiveSheet
$win = $excel.ActiveWindow

# --- sheetView: zoom + frozen-pane anchor + active selection ---------------
$win.Zoom = 90

# The pane is already frozen at row 1 (ySplit=1); re-point its scrolled
# top-left corner back to A2 instead of A76.
$win.ScrollRow = 2
$win.ScrollColumn = 1

# Active cell/selection in the (frozen-below) pane moves to E15.
$ws.Range("E15").Select()

# --- column widths -----------------------------------------------------------
# Former single run: <col min="4" max="9" width="20.1796875"/>
# now split into per-column widths; col 10 (J) also gets a new width.
$ws.Columns("D").ColumnWidth = 4.16666666666667
$ws.Columns("E").ColumnWidth = 57
$ws.Columns("F").ColumnWidth = 4.33333333333333
$ws.Columns("G").ColumnWidth = 3.66666666666667
$ws.Columns("H").ColumnWidth = 57
$ws.Columns("I").ColumnWidth = 3.33333333333333
$ws.Columns("J").ColumnWidth = 29.3333333333333
